$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.144.64"
$ws.Range("E2").Value = "  +2.66%  "
$ws.Range("D3").Value = "2.358.04"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'312.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").Value = "'107.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "'0.610"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.85%  "
$ws.Range("D10").Value = "'40.88"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("D11").Value = "'0.0917"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "'8.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.74%  "
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "'0.973"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.03%  "
$ws.Range("D15").Value = "2.715.27"
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("D17").Value = "2.358.40"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "45.091.87"
$ws.Range("E18").Value = "  +2.76%  "
$ws.Range("D19").Value = "'14.33"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +10.64%  "
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").Value = "'7.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.99%  "
$ws.Range("D22").Value = "'73.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.78%  "
$ws.Range("D23").Value = "'3.50"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").Value = "'259.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.30%  "
$ws.Range("E25").Value = "  +1.40%  "
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").Value = "'11.04"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("D28").Value = "'7.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.55%  "
$ws.Range("E29").Value = "  +0.27%  "
$ws.Range("D30").Value = "'0.0968"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.79%  "
$ws.Range("D31").Value = "'22.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("D32").Value = "'37.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.94%  "
$ws.Range("D33").Value = "'167.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("D34").Value = "'2.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.64%  "
$ws.Range("E35").Value = "  -1.70%  "
$ws.Range("D36").Value = "'0.116"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("E37").Value = "  -1.32%  "
$ws.Range("D38").Value = "'3.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.27%  "
$ws.Range("E39").Value = "  -3.14%  "
$ws.Range("D40").Value = "'2.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.89%  "
$ws.Range("D41").Value = "'1.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.72%  "
$ws.Range("D42").Value = "'99.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.86%  "
$ws.Range("D43").Value = "'69.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.17%  "
$ws.Range("D44").Value = "'0.228"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.86%  "
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").Value = "'12.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.22%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.35%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.842.93"
$ws.Range("E47").Value = "  +10.99%  "
$ws.Range("D48").Value = "'83.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.89%  "
$ws.Range("D49").Value = "'5.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.60%  "
$ws.Range("D50").Value = "'110.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.12%  "
$ws.Range("D51").Value = "'9.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.87%  "
